$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4584.857
$ws.Range("I64").Value = 4765.8335
$ws.Range("J64").Value = 3499
$ws.Range("K64").Value = 4765.8335
$ws.Range("L64").Value = 3499
$ws.Range("M64").Value = -4517.8335
$ws.Range("N64").Value = -3995
$ws.Range("H67").Value = 4584.857
$ws.Range("I67").Value = 4765.8335
$ws.Range("J67").Value = 3499
$ws.Range("K67").Value = 4765.8335
$ws.Range("L67").Value = 3499
$ws.Range("M67").Value = -3907.8335
$ws.Range("N67").Value = -5215
$ws.Range("H69").Value = 18625
$ws.Range("I69").Value = 19661.666
$ws.Range("K69").Value = 58984.99800000001
$ws.Range("M69").Value = -58110.99800000001
$ws.Range("H72").Value = 18625
$ws.Range("I72").Value = 19661.666
$ws.Range("K72").Value = 176954.994
$ws.Range("M72").Value = -172586.994
$ws.Range("H113").Value = 9579.799999999999
$ws.Range("I113").Value = 8999.5
$ws.Range("J113").Value = 9966.666999999999
$ws.Range("K113").Value = 8999.5
$ws.Range("L113").Value = 9966.666999999999
$ws.Range("M113").Value = -5745.5
$ws.Range("N113").Value = -16474.667
$ws.Range("H138").Value = 8235.436
$ws.Range("I138").Value = 6896.357
$ws.Range("J138").Value = 8985.32
$ws.Range("K138").Value = 20689.071
$ws.Range("L138").Value = 26955.96
$ws.Range("M138").Value = -15549.071
$ws.Range("N138").Value = -37235.96

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16863
$ws.Range("I32").Value = 16863
$ws.Range("K32").Value = 16863
$ws.Range("M32").Value = -16576
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 781.5484
$ws.Range("I97").Value = 704.0769
$ws.Range("J97").Value = 1184.4
$ws.Range("K97").Value = 704.0769
$ws.Range("L97").Value = 1184.4
$ws.Range("M97").Value = -208.0769
$ws.Range("N97").Value = -2176.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3068.125
$ws.Range("J31").Value = 3272.5
$ws.Range("L31").Value = 3272.5
$ws.Range("N31").Value = -3862.5
$ws.Range("H34").Value = 3068.125
$ws.Range("J34").Value = 3272.5
$ws.Range("L34").Value = 3272.5
$ws.Range("N34").Value = -3676.5
$ws.Range("H62").Value = 5584.6665
$ws.Range("J62").Value = 5584.6665
$ws.Range("L62").Value = 5584.6665
$ws.Range("N62").Value = -6832.6665
$ws.Range("H65").Value = 5584.6665
$ws.Range("J65").Value = 5584.6665
$ws.Range("L65").Value = 27923.3325
$ws.Range("N65").Value = -34163.3325
$ws.Range("H132").Value = 13344019
$ws.Range("I132").Value = 16679141
$ws.Range("J132").Value = 3532.6667
$ws.Range("K132").Value = 50037423
$ws.Range("L132").Value = 10598.0001
$ws.Range("M132").Value = -50034893
$ws.Range("N132").Value = -15658.0001
$ws.Range("H141").Value = 1074824.2
$ws.Range("J141").Value = 1074824.2
$ws.Range("L141").Value = 1074824.2
$ws.Range("N141").Value = -1085184.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3892.3333
$ws.Range("I3").Value = 3892.3333
$ws.Range("K3").Value = 11676.9999
$ws.Range("M3").Value = -11564.9999
$ws.Range("H5").Value = 103579.1
$ws.Range("I5").Value = 4223.25
$ws.Range("J5").Value = 501002.5
$ws.Range("K5").Value = 12669.75
$ws.Range("L5").Value = 1503007.5
$ws.Range("M5").Value = -12557.75
$ws.Range("N5").Value = -1503231.5
$ws.Range("H17").Value = 596.3333
$ws.Range("I17").Value = 192.66667
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 578.00001
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -409.00001
$ws.Range("N17").Value = -3338
$ws.Range("H52").Value = 7444.5
$ws.Range("J52").Value = 7444.5
$ws.Range("L52").Value = 22333.5
$ws.Range("N52").Value = -22865.5
$ws.Range("H86").Value = 161.5
$ws.Range("I86").Value = 161.5
$ws.Range("K86").Value = 484.5
$ws.Range("M86").Value = 701.5
$ws.Range("H89").Value = 161.5
$ws.Range("I89").Value = 161.5
$ws.Range("K89").Value = 1453.5
$ws.Range("M89").Value = 4474.5
$ws.Range("H92").Value = 366.5
$ws.Range("I92").Value = 375
$ws.Range("K92").Value = 1125
$ws.Range("M92").Value = 123
$ws.Range("H114").Value = 1870.9
$ws.Range("I114").Value = 1150
$ws.Range("J114").Value = 2051.125
$ws.Range("K114").Value = 3450
$ws.Range("L114").Value = 6153.375
$ws.Range("M114").Value = -196
$ws.Range("N114").Value = -12661.375
$ws.Range("H122").Value = 15347.571
$ws.Range("I122").Value = 17590.666
$ws.Range("K122").Value = 158315.994
$ws.Range("M122").Value = -155865.994
$ws.Range("H135").Value = 103579.1
$ws.Range("I135").Value = 4223.25
$ws.Range("J135").Value = 501002.5
$ws.Range("K135").Value = 38009.25
$ws.Range("L135").Value = 4509022.5
$ws.Range("M135").Value = -35474.25
$ws.Range("N135").Value = -4514092.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7054.091
$ws.Range("I80").Value = 3110.5
$ws.Range("J80").Value = 7930.4443
$ws.Range("K80").Value = 3110.5
$ws.Range("L80").Value = 7930.4443
$ws.Range("M80").Value = -2112.5
$ws.Range("N80").Value = -9926.444299999999
$ws.Range("H83").Value = 7054.091
$ws.Range("I83").Value = 3110.5
$ws.Range("J83").Value = 7930.4443
$ws.Range("K83").Value = 15552.5
$ws.Range("L83").Value = 39652.2215
$ws.Range("M83").Value = -10560.5
$ws.Range("N83").Value = -49636.2215

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9238
$ws.Range("I61").Value = 8391.571
$ws.Range("K61").Value = 8391.571
$ws.Range("M61").Value = -8189.571
$ws.Range("H93").Value = 1273.0769
$ws.Range("I93").Value = 1212.5834
$ws.Range("J93").Value = 1999
$ws.Range("K93").Value = 1212.5834
$ws.Range("L93").Value = 1999
$ws.Range("M93").Value = 35.41660000000002
$ws.Range("N93").Value = -4495
$ws.Range("H113").Value = 9238
$ws.Range("I113").Value = 8391.571
$ws.Range("K113").Value = 8391.571
$ws.Range("M113").Value = -6221.571
